$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 46, shifting existing rows 46-113 down to 47-114
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with data
$ws.Cells.Item(46, 1).Value = 5
$ws.Cells.Item(46, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(46, 3).Value = "Maule"
$ws.Cells.Item(46, 4).Value = 45210
$ws.Cells.Item(46, 5).Value = 7
$ws.Cells.Item(46, 6).Value = 300000000
$ws.Cells.Item(46, 7).Value = "Espárragos"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 4000
$ws.Cells.Item(46, 11).Value = 1100
$ws.Cells.Item(46, 12).Value = 1200
$ws.Cells.Item(46, 13).Value = 1150
$ws.Cells.Item(46, 14).Value = "`$/kilo"
$ws.Cells.Item(46, 15).Value = "Provincia de Linares"
$ws.Cells.Item(46, 16).Value = 1150
$ws.Cells.Item(46, 17).Value = 1
$ws.Cells.Item(46, 18).Value = "Hortaliza"
